{"js": "// Replace the 100 arithmetic-expression cells in the single table with\n// their updated expressions. The mapping below is positional (row-major,\n// left-to-right, top-to-bottom, matching the table's natural cell order)\n// because a few original expressions repeat (e.g. \"71-67=\") but must be\n// replaced with different new values depending on which occurrence it is.\nconst replacements = [\n  [\"8+33=\", \"82-13=\"], [\"1+21=\", \"56-12=\"], [\"77-25=\", \"60-27=\"], [\"57-57=\", \"88-8=\"], [\"85-3=\", \"72+10=\"],\n  [\"52+31=\", \"92-53=\"], [\"0+64=\", \"62+28=\"], [\"80-56=\", \"46+3=\"], [\"76-41=\", \"87-4=\"], [\"39+60=\", \"69+16=\"],\n  [\"44-7=\", \"25+30=\"], [\"62-10=\", \"81-39=\"], [\"75-18=\", \"8+77=\"], [\"28+24=\", \"99-75=\"], [\"89-83=\", \"5+46=\"],\n  [\"1+67=\", \"30+21=\"], [\"3+10=\", \"75-55=\"], [\"69-32=\", \"99-41=\"], [\"95-94=\", \"29+44=\"], [\"93-55=\", \"47+8=\"],\n  [\"3+83=\", \"9+90=\"], [\"97-43=\", \"54-47=\"], [\"34+6=\", \"83-45=\"], [\"8+64=\", \"36-18=\"], [\"6+52=\", \"95-93=\"],\n  [\"74+12=\", \"30+9=\"], [\"92-51=\", \"93-64=\"], [\"13+53=\", \"26+14=\"], [\"25-6=\", \"61-46=\"], [\"70-48=\", \"21-7=\"],\n  [\"25+13=\", \"33+63=\"], [\"71-67=\", \"60-40=\"], [\"84-37=\", \"26-21=\"], [\"29-28=\", \"82-41=\"], [\"98-23=\", \"60+27=\"],\n  [\"12+29=\", \"65-3=\"], [\"39+32=\", \"34-33=\"], [\"11+59=\", \"72+26=\"], [\"12+84=\", \"84-56=\"], [\"85-38=\", \"76-1=\"],\n  [\"71-53=\", \"38-30=\"], [\"32+20=\", \"25+17=\"], [\"66+24=\", \"50-42=\"], [\"64-40=\", \"80-61=\"], [\"55+42=\", \"59-15=\"],\n  [\"40+10=\", \"70+20=\"], [\"54+24=\", \"54+10=\"], [\"65-57=\", \"21+1=\"], [\"15+71=\", \"19+70=\"], [\"9+47=\", \"77-53=\"],\n  [\"76-53=\", \"9+8=\"], [\"80-16=\", \"76-18=\"], [\"40+54=\", \"85-32=\"], [\"87+4=\", \"78-25=\"], [\"78-9=\", \"98-50=\"],\n  [\"72-62=\", \"15+40=\"], [\"39-0=\", \"47+42=\"], [\"90-4=\", \"13+31=\"], [\"34+18=\", \"80-27=\"], [\"71+19=\", \"17+33=\"],\n  [\"75-69=\", \"9+69=\"], [\"82-12=\", \"82-46=\"], [\"32-29=\", \"47+28=\"], [\"71-67=\", \"24-4=\"], [\"68-29=\", \"35-5=\"],\n  [\"91-69=\", \"13+23=\"], [\"93-9=\", \"62-56=\"], [\"64+18=\", \"22+66=\"], [\"36-17=\", \"14+21=\"], [\"39+7=\", \"31+27=\"],\n  [\"55-24=\", \"30+8=\"], [\"44+4=\", \"57-18=\"], [\"80-1=\", \"82-54=\"], [\"11+5=\", \"32+57=\"], [\"58+32=\", \"77-7=\"],\n  [\"82-35=\", \"47+40=\"], [\"81-9=\", \"75-27=\"], [\"44+52=\", \"28+30=\"], [\"11+13=\", \"73-19=\"], [\"98-2=\", \"86-31=\"],\n  [\"78-56=\", \"53-8=\"], [\"44+2=\", \"39+20=\"], [\"66-37=\", \"56-4=\"], [\"2+38=\", \"31+16=\"], [\"97-93=\", \"13+83=\"],\n  [\"49-3=\", \"20-19=\"], [\"93-40=\", \"93-34=\"], [\"63+17=\", \"77+10=\"], [\"1+63=\", \"74-71=\"], [\"71-6=\", \"44-13=\"],\n  [\"54-0=\", \"44+35=\"], [\"85-70=\", \"52+13=\"], [\"47+16=\", \"29+12=\"], [\"27+36=\", \"2+70=\"], [\"4+17=\", \"75-42=\"],\n  [\"36+24=\", \"42+12=\"], [\"32-17=\", \"72-50=\"], [\"40+40=\", \"51-41=\"], [\"25+74=\", \"38-10=\"], [\"99-72=\", \"52-37=\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst oldValues = table.values; // 2-D array, row-major, matches cell order\nlet k = 0;\nconst newValues = oldValues.map((row) =>\n  row.map((cellText) => {\n    const [expectedOld, newText] = replacements[k];\n    if (cellText !== expectedOld) {\n      throw new Error(\n        `Cell #${k} text \"${cellText}\" does not match expected \"${expectedOld}\"`\n      );\n    }\n    k += 1;\n    return newText;\n  })\n);\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-expression cells in the single table with\n# their updated expressions. The mapping below is positional (row-major,\n# left-to-right, top-to-bottom, matching Cell(row, col) order) because a\n# few original expressions repeat (e.g. \"71-67=\") but must be replaced\n# with different new values depending on which occurrence it is.\n\n$replacements = @(\n  @(\"8+33=\", \"82-13=\"), @(\"1+21=\", \"56-12=\"), @(\"77-25=\", \"60-27=\"), @(\"57-57=\", \"88-8=\"), @(\"85-3=\", \"72+10=\"),\n  @(\"52+31=\", \"92-53=\"), @(\"0+64=\", \"62+28=\"), @(\"80-56=\", \"46+3=\"), @(\"76-41=\", \"87-4=\"), @(\"39+60=\", \"69+16=\"),\n  @(\"44-7=\", \"25+30=\"), @(\"62-10=\", \"81-39=\"), @(\"75-18=\", \"8+77=\"), @(\"28+24=\", \"99-75=\"), @(\"89-83=\", \"5+46=\"),\n  @(\"1+67=\", \"30+21=\"), @(\"3+10=\", \"75-55=\"), @(\"69-32=\", \"99-41=\"), @(\"95-94=\", \"29+44=\"), @(\"93-55=\", \"47+8=\"),\n  @(\"3+83=\", \"9+90=\"), @(\"97-43=\", \"54-47=\"), @(\"34+6=\", \"83-45=\"), @(\"8+64=\", \"36-18=\"), @(\"6+52=\", \"95-93=\"),\n  @(\"74+12=\", \"30+9=\"), @(\"92-51=\", \"93-64=\"), @(\"13+53=\", \"26+14=\"), @(\"25-6=\", \"61-46=\"), @(\"70-48=\", \"21-7=\"),\n  @(\"25+13=\", \"33+63=\"), @(\"71-67=\", \"60-40=\"), @(\"84-37=\", \"26-21=\"), @(\"29-28=\", \"82-41=\"), @(\"98-23=\", \"60+27=\"),\n  @(\"12+29=\", \"65-3=\"), @(\"39+32=\", \"34-33=\"), @(\"11+59=\", \"72+26=\"), @(\"12+84=\", \"84-56=\"), @(\"85-38=\", \"76-1=\"),\n  @(\"71-53=\", \"38-30=\"), @(\"32+20=\", \"25+17=\"), @(\"66+24=\", \"50-42=\"), @(\"64-40=\", \"80-61=\"), @(\"55+42=\", \"59-15=\"),\n  @(\"40+10=\", \"70+20=\"), @(\"54+24=\", \"54+10=\"), @(\"65-57=\", \"21+1=\"), @(\"15+71=\", \"19+70=\"), @(\"9+47=\", \"77-53=\"),\n  @(\"76-53=\", \"9+8=\"), @(\"80-16=\", \"76-18=\"), @(\"40+54=\", \"85-32=\"), @(\"87+4=\", \"78-25=\"), @(\"78-9=\", \"98-50=\"),\n  @(\"72-62=\", \"15+40=\"), @(\"39-0=\", \"47+42=\"), @(\"90-4=\", \"13+31=\"), @(\"34+18=\", \"80-27=\"), @(\"71+19=\", \"17+33=\"),\n  @(\"75-69=\", \"9+69=\"), @(\"82-12=\", \"82-46=\"), @(\"32-29=\", \"47+28=\"), @(\"71-67=\", \"24-4=\"), @(\"68-29=\", \"35-5=\"),\n  @(\"91-69=\", \"13+23=\"), @(\"93-9=\", \"62-56=\"), @(\"64+18=\", \"22+66=\"), @(\"36-17=\", \"14+21=\"), @(\"39+7=\", \"31+27=\"),\n  @(\"55-24=\", \"30+8=\"), @(\"44+4=\", \"57-18=\"), @(\"80-1=\", \"82-54=\"), @(\"11+5=\", \"32+57=\"), @(\"58+32=\", \"77-7=\"),\n  @(\"82-35=\", \"47+40=\"), @(\"81-9=\", \"75-27=\"), @(\"44+52=\", \"28+30=\"), @(\"11+13=\", \"73-19=\"), @(\"98-2=\", \"86-31=\"),\n  @(\"78-56=\", \"53-8=\"), @(\"44+2=\", \"39+20=\"), @(\"66-37=\", \"56-4=\"), @(\"2+38=\", \"31+16=\"), @(\"97-93=\", \"13+83=\"),\n  @(\"49-3=\", \"20-19=\"), @(\"93-40=\", \"93-34=\"), @(\"63+17=\", \"77+10=\"), @(\"1+63=\", \"74-71=\"), @(\"71-6=\", \"44-13=\"),\n  @(\"54-0=\", \"44+35=\"), @(\"85-70=\", \"52+13=\"), @(\"47+16=\", \"29+12=\"), @(\"27+36=\", \"2+70=\"), @(\"4+17=\", \"75-42=\"),\n  @(\"36+24=\", \"42+12=\"), @(\"32-17=\", \"72-50=\"), @(\"40+40=\", \"51-41=\"), @(\"25+74=\", \"38-10=\"), @(\"99-72=\", \"52-37=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n$k = 0\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cellRange = $cell.Range\n    # Cell.Range.Text includes a trailing cell-mark; strip it for comparison.\n    $current = $cellRange.Text.TrimEnd([char]7, [char]13)\n\n    $expectedOld = $replacements[$k][0]\n    $newText = $replacements[$k][1]\n\n    if ($current -ne $expectedOld) {\n      throw \"Cell ($r,$c) text '$current' does not match expected '$expectedOld'\"\n    }\n\n    $cellRange.Text = $newText\n    $k++\n  }\n}\n"}
